$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append "mm" unit suffix to the dimension values in columns B and C (rows 2-9),
# converting them from numbers to text.
for ($r = 2; $r -le 9; $r++) {
    $bCell = $ws.Cells.Item($r, 2)
    $bCell.Value = ([string]($bCell.Value2)) + "mm"
}

for ($r = 2; $r -le 9; $r++) {
    $cCell = $ws.Cells.Item($r, 3)
    $cCell.Value = ([string]($cCell.Value2)) + "mm"
}
